# February 24, 2021 update -- append 13 new daily rows (336-348) for Germany
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A-C: constant identifier columns (iso_code/continent/location) for every new row
$ws.Range("A336:A348").Value = "DEU"
$ws.Range("B336:B348").Value = "Europe"
$ws.Range("C336:C348").Value = "Germany"

# D: date column -- extends the existing "previous day + 1" fill pattern
$ws.Range("D336:D348").Formula = "=D335+1"

# F/H: raw daily inputs (new cases, total tests); K-P stay 0 for this sheet
$ws.Range("F336").Value = 9860
$ws.Range("H336").Value = 64191
$ws.Range("K336:P336").Value = 0
$ws.Range("F337").Value = 8354
$ws.Range("H337").Value = 64672
$ws.Range("K337:P337").Value = 0
$ws.Range("F338").Value = 6114
$ws.Range("H338").Value = 64960
$ws.Range("K338:P338").Value = 0
$ws.Range("F339").Value = 4426
$ws.Range("H339").Value = 65076
$ws.Range("K339:P339").Value = 0
$ws.Range("F340").Value = 3856
$ws.Range("H340").Value = 65604
$ws.Range("K340:P340").Value = 0
$ws.Range("F341").Value = 7556
$ws.Range("H341").Value = 66164
$ws.Range("K341:P341").Value = 0
$ws.Range("F342").Value = 10207
$ws.Range("H342").Value = 66698
$ws.Range("K342:P342").Value = 0
$ws.Range("F343").Value = 9113
$ws.Range("H343").Value = 67206
$ws.Range("K343:P343").Value = 0
$ws.Range("F344").Value = 9164
$ws.Range("H344").Value = 67696
$ws.Range("K344:P344").Value = 0
$ws.Range("F345").Value = 7676
$ws.Range("H345").Value = 67841
$ws.Range("K345:P345").Value = 0
$ws.Range("F346").Value = 4369
$ws.Range("H346").Value = 67903
$ws.Range("K346:P346").Value = 0
$ws.Range("F347").Value = 3883
$ws.Range("H347").Value = 68318
$ws.Range("K347:P347").Value = 0
$ws.Range("F348").Value = 8007
$ws.Range("H348").Value = 68740
$ws.Range("K348:P348").Value = 0

# E: running cumulative total = previous cumulative total + new F value
$ws.Range("E336").Formula = "=E335+F336"
$ws.Range("E337").Formula = "=E336+F337"
$ws.Range("E338").Formula = "=E337+F338"
$ws.Range("E339").Formula = "=E338+F339"
$ws.Range("E340").Formula = "=E339+F340"
$ws.Range("E341").Formula = "=E340+F341"
$ws.Range("E342").Formula = "=E341+F342"
$ws.Range("E343").Formula = "=E342+F343"
$ws.Range("E344").Formula = "=E343+F344"
$ws.Range("E345").Formula = "=E344+F345"
$ws.Range("E346").Formula = "=E345+F346"
$ws.Range("E347").Formula = "=E346+F347"
$ws.Range("E348").Formula = "=E347+F348"

# G: trailing 7-day average of F
$ws.Range("G336").Formula = "=SUM(F330:F336)/7"
$ws.Range("G337").Formula = "=SUM(F331:F337)/7"
$ws.Range("G338").Formula = "=SUM(F332:F338)/7"
$ws.Range("G339").Formula = "=SUM(F333:F339)/7"
$ws.Range("G340").Formula = "=SUM(F334:F340)/7"
$ws.Range("G341").Formula = "=SUM(F335:F341)/7"
$ws.Range("G342").Formula = "=SUM(F336:F342)/7"
$ws.Range("G343").Formula = "=SUM(F337:F343)/7"
$ws.Range("G344").Formula = "=SUM(F338:F344)/7"
$ws.Range("G345").Formula = "=SUM(F339:F345)/7"
$ws.Range("G346").Formula = "=SUM(F340:F346)/7"
$ws.Range("G347").Formula = "=SUM(F341:F347)/7"
$ws.Range("G348").Formula = "=SUM(F342:F348)/7"

# I: day-over-day change in H
$ws.Range("I336").Formula = "=H336-H335"
$ws.Range("I337").Formula = "=H337-H336"
$ws.Range("I338").Formula = "=H338-H337"
$ws.Range("I339").Formula = "=H339-H338"
$ws.Range("I340").Formula = "=H340-H339"
$ws.Range("I341").Formula = "=H341-H340"
$ws.Range("I342").Formula = "=H342-H341"
$ws.Range("I343").Formula = "=H343-H342"
$ws.Range("I344").Formula = "=H344-H343"
$ws.Range("I345").Formula = "=H345-H344"
$ws.Range("I346").Formula = "=H346-H345"
$ws.Range("I347").Formula = "=H347-H346"
$ws.Range("I348").Formula = "=H348-H347"

# J: trailing 7-day average of I
$ws.Range("J336").Formula = "=SUM(I330:I336)/7"
$ws.Range("J337").Formula = "=SUM(I331:I337)/7"
$ws.Range("J338").Formula = "=SUM(I332:I338)/7"
$ws.Range("J339").Formula = "=SUM(I333:I339)/7"
$ws.Range("J340").Formula = "=SUM(I334:I340)/7"
$ws.Range("J341").Formula = "=SUM(I335:I341)/7"
$ws.Range("J342").Formula = "=SUM(I336:I342)/7"
$ws.Range("J343").Formula = "=SUM(I337:I343)/7"
$ws.Range("J344").Formula = "=SUM(I338:I344)/7"
$ws.Range("J345").Formula = "=SUM(I339:I345)/7"
$ws.Range("J346").Formula = "=SUM(I340:I346)/7"
$ws.Range("J347").Formula = "=SUM(I341:I347)/7"
$ws.Range("J348").Formula = "=SUM(I342:I348)/7"

# Leave the sheet scrolled/selected where the author left it when saving
$ws.Range("F342").Select()
